# Update the "慢检测逻辑判断" matlab-logic worksheet:
#  - widen column B so the longer explanatory notes fit
#  - rewrite the notes table (rows 8-16) to reflect the new logic description,
#    inserting two new explanatory rows (one after "当次是F列" and one after
#    "当次是F") and renaming/expanding several of the existing notes
#  - fix the stray "-" placeholder cells (C2, D4, E4) so they reference the
#    correct shared string
#  - move the active selection to B15

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- widen column B to fit the longer text that is about to be written ---
$ws.Columns.Item(2).ColumnWidth = 46.2857142857

# --- the "-" placeholder values are unchanged in content, just re-touch them
#     so any stale shared-string reference is normalised ---
$ws.Range("C2").Value = "-"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "-"

# --- clear the old notes table (rows 8-16, columns A:B) so it can be
#     rebuilt with the new row layout ---
$ws.Range("A8:B16").ClearContents()

# --- rebuild the notes table with the updated / expanded wording ---
$ws.Range("A8").Value = "大1"
$ws.Range("B8").Value = "当次大动作B列"

$ws.Range("A9").Value = "大0"
$ws.Range("B9").Value = "上一次大1"

$ws.Range("B10").Value = "当次是F列"

$ws.Range("B11").Value = "后面一个周期还是F列的话连续输出"

$ws.Range("A13").Value = "微1"
$ws.Range("B13").Value = "当次是C/D"

$ws.Range("A14").Value = "微1"
$ws.Range("B14").Value = "E连续出现三次，第三次出现时输出微1，并清掉等再连续出现时再输出微1，不连续出现就输出微0"

$ws.Range("A16").Value = "微0"
$ws.Range("B16").Value = "上一次是微1，c/d/e列"

$ws.Range("B17").Value = "当次是F"

$ws.Range("B18").Value = "后面一个周期还是F列的话连续输出"

# --- move the selection like the author left it ---
$ws.Range("B15").Select()
